$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.945.51"
$ws.Range("E2").Value = "  -3.19%  "

$ws.Range("D3").Value = "2.993.13"
$ws.Range("E3").Value = "  -3.46%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.14%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "3.007.88"
$ws.Range("E9").Value = "  -3.22%  "

$ws.Range("E10").Value = "  -1.93%  "

$ws.Range("E11").Value = "  -6.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.369"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.91%  "

$ws.Range("D13").Value = "3.519.88"
$ws.Range("E13").Value = "  -3.63%  "

$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").Value = "61.965.63"
$ws.Range("E15").Value = "  -3.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.16%  "

$ws.Range("D17").Value = "3.003.04"
$ws.Range("E17").Value = "  -3.66%  "

$ws.Range("E18").Value = "  -3.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.75%  "

$ws.Range("D25").Value = "3.118.46"
$ws.Range("E25").Value = "  -4.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.469"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.45%  "

$ws.Range("E27").Value = "  -1.73%  "

$ws.Range("E28").Value = "  -0.58%  "

$ws.Range("D29").Value = "0.0₃0936"
$ws.Range("E29").Value = "  -7.13%  "

$ws.Range("E30").Value = "  -7.56%  "

$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("E32").Value = "  -3.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "161.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.84%  "

$ws.Range("E35").Value = "  -4.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.87%  "

$ws.Range("E37").Value = "  -3.99%  "

$ws.Range("E38").Value = "  -4.61%  "

$ws.Range("E39").Value = "  -5.09%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.83%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.421.24"
$ws.Range("E42").Value = "  -6.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.673"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.81%  "

$ws.Range("E45").Value = "  -2.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("E48").Value = "  -3.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "269.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0953"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.18%  "
